$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number-format on the numeric-looking Price column cells we touch,
# so Excel stores them as text (matching the original inline-string cells)
# instead of silently re-parsing them as numbers and dropping trailing zeros.
$ws.Range('D2').Value = '97.164.56'
$ws.Range('E2').Value = '  +0.68%  '
$ws.Range('D3').Value = '3.706.16'
$ws.Range('E3').Value = '  +0.88%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '238.04'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -1.12%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.90'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +3.62%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '661.14'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -0.22%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.425'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +1.00%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.999'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -0.06%  '
$ws.Range('E10').Value = '  -0.67%  '
$ws.Range('D11').Value = '3.703.84'
$ws.Range('E11').Value = '  +0.98%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000323'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +20.58%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '44.30'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -2.16%  '
$ws.Range('E14').Value = '  +1.65%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.82'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -0.34%  '
$ws.Range('D16').Value = '4.396.39'
$ws.Range('E16').Value = '  +0.86%  '
$ws.Range('D17').Value = '96.968.61'
$ws.Range('E17').Value = '  +0.75%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '9.10'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +2.93%  '
$ws.Range('D19').Value = '3.709.60'
$ws.Range('E19').Value = '  +0.85%  '
$ws.Range('E20').Value = '  +1.64%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '18.70'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +1.42%  '
$ws.Range('E22').Value = '  -3.67%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '522.20'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -0.40%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.44'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +0.36%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.0000218'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +7.05%  '
$ws.Range('E26').Value = '  -0.98%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '101.85'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +0.09%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.192'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +13.60%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '13.60'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +4.38%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '12.83'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +2.90%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.05'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -0.04%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.999'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -0.03%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.88'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -1.27%  '
$ws.Range('E35').Value = '  +0.25%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '32.25'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -1.41%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '653.66'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +1.46%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.594'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +0.55%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '8.88'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +1.53%  '
$ws.Range('E40').Value = '  +0.06%  '
$ws.Range('E41').Value = '  +2.67%  '
$ws.Range('B42').Value = 'Filecoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.79'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +6.10%  '
$ws.Range('B43').Value = 'ImmutableX'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.04'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +3.86%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '40.42'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -9.08%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.482'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +5.68%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.970'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +0.45%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0459'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -0.79%  '
$ws.Range('E48').Value = '  +0.73%  '
$ws.Range('E49').Value = '  -0.14%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.71'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +1.45%  '
$ws.Range('E51').Value = '  -3.33%  '
